$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Full replacement data for rows 2-7 (columns A-T).
# Row 2 and 3 come from splitting the old "ECs" row (target cluster ECs vs sCs).
# Row 4 and 5 come from splitting the old "FAPs" row (target cluster ECs vs sCs).
# Row 6 and 7 come from splitting the old "sCs" row (target cluster ECs vs sCs).
$data = @(
    @("ECs",  "Gnai2", "P2ry12", "ECs", 3, 1, 184.0626906666667,  552.188072,         0.6510505751503485, 0.6510505751503486, 2, 0.6666666666666666, 6.66513,            19.99539,  0.9556042431346229, 0.9556042431346229, 1226.80176144312,  11041.21585298808, 0.6221466921089097, 0.6221466921089098),
    @("ECs",  "Gnai2", "P2ry12", "sCs", 3, 1, 184.0626906666667,  552.188072,         0.6510505751503485, 0.6510505751503486, 3, 1,                  0.3096506666666667, 0.928952,  0.04439575686537718, 0.04439575686537718, 56.99513487339378, 512.956213860544,  0.02890388304143885, 0.02890388304143885),
    @("FAPs", "Gnai2", "P2ry12", "ECs", 3, 1, 57.4434,            172.3302,           0.2031838091312023, 0.2031838091312023, 2, 0.6666666666666666, 6.66513,            19.99539,  0.9556042431346229, 0.9556042431346229, 382.867728642,     3445.809557778,    0.1941633101420323, 0.1941633101420323),
    @("FAPs", "Gnai2", "P2ry12", "sCs", 3, 1, 57.4434,            172.3302,           0.2031838091312023, 0.2031838091312023, 3, 1,                  0.3096506666666667, 0.928952,  0.04439575686537718, 0.04439575686537718, 17.7873871056,     160.0864839504,    0.00902049898917006, 0.00902049898917006),
    @("sCs",  "Gnai2", "P2ry12", "ECs", 3, 1, 41.21033366666666,  123.631001,         0.1457656157184491, 0.1457656157184491, 2, 0.6666666666666666, 6.66513,            19.99539,  0.9556042431346229, 0.9556042431346229, 274.67223123171,   2472.05008108539,  0.1392942408836809, 0.1392942408836809),
    @("sCs",  "Gnai2", "P2ry12", "sCs", 3, 1, 41.21033366666666,  123.631001,         0.1457656157184491, 0.1457656157184491, 3, 1,                  0.3096506666666667, 0.928952,  0.04439575686537718, 0.04439575686537718, 12.76080729343911, 114.847265640952,  0.006471374834768269, 0.006471374834768269)
)

for ($i = 0; $i -lt $data.Length; $i++) {
    $r = $i + 2
    $rowValues = $data[$i]
    for ($c = 1; $c -le $rowValues.Length; $c++) {
        $ws.Cells.Item($r, $c).Value = $rowValues[$c - 1]
    }
}
